# Fixed Organic N extraction from excel table (id 20 instead of 15).
# Updated the database values and reset the active view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sample_ID value changed
$ws.Range("B4").Value = 104

# Parameter values pulled fresh from the (fixed) database query
$ws.Range("B18").Value = 17   # N_total
$ws.Range("B19").Value = 8    # N_NH4
$ws.Range("B20").Value = 9    # N_NO3
$ws.Range("B21").Value = 2    # N_Organic
$ws.Range("B22").Value = 8    # P
$ws.Range("B23").Value = 45   # K
$ws.Range("B24").Value = 12   # Ca
$ws.Range("B25").Value = 15   # Mg
$ws.Range("B26").Value = 42   # S
$ws.Range("B27").Value = 2    # Fe
$ws.Range("B28").Value = 4    # Mn
$ws.Range("B29").Value = 1    # B
$ws.Range("B32").Value = 1    # Mo

# Reset scroll position to the top of the sheet and move the selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
